$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the old "aithyia" / "apotheosis/heroization" header values and
# rebuild row 1 headers (shifting columns B..J right by one and updating J1)
$ws.Range("B1").Value = "web-footed"
$ws.Range("C1").Value = "white"
$ws.Range("D1").Value = "black"
$ws.Range("E1").Value = "bird of prey"
$ws.Range("F1").Value = "duck"
$ws.Range("G1").Value = "woman"
$ws.Range("H1").Value = "man"
$ws.Range("I1").Value = "dive into the sea"
$ws.Range("J1").Value = "myth"

# Row 2: aithyia (gull)
$ws.Range("A2").Value = "aithyia (gull)"
$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = "unknown"
$ws.Range("G2").Value = "yes"
$ws.Range("I2").Value = "yes"
$ws.Range("J2").Value = "Ino"

# Row 3: memnon (ruff)
$ws.Range("A3").Value = "memnon (ruff)"
$ws.Range("H3").Value = "yes"
$ws.Range("J3").Value = "Companions of Memnon"

# Row 4: ortyx (quail)
$ws.Range("A4").Value = "ortyx (quail)"
$ws.Range("G4").Value = "yes"
$ws.Range("I4").Value = "yes"
$ws.Range("J4").Value = "Asteria"

# Row 5: erodioi(herons)
$ws.Range("A5").Value = "erodioi(herons)"
$ws.Range("C5").Value = "yes"
$ws.Range("H5").Value = "yes"
$ws.Range("J5").Value = "Companions of Diomedes"

# Row 6: erodioi(shearwaters)
$ws.Range("A6").Value = "erodioi(shearwaters)"
$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = "yes"
$ws.Range("H6").Value = "yes"
$ws.Range("J6").Value = "Companions of Diomedes (after recant in 1918)"

# Column A width change (target stored width 18.7109375; engine quantizes
# column widths to 1/6-character steps, so 17.8 is the closest input that
# rounds to the nearest achievable stored width of 18.666666666666668)
$ws.Columns.Item(1).ColumnWidth = 17.8

$ws.Range("A7").Select()
